$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 147, pushing the existing data
# (rows 147-260) down to rows 148-261.
$ws.Rows("147").Insert()

# Populate the newly inserted row with the new price-record data.
$ws.Range("A147").Value = 3
$ws.Range("B147").Value = "Femacal de La Calera"
$ws.Range("C147").Value = "Coquimbo"
$ws.Range("D147").Value = 44981
$ws.Range("E147").Value = 5
$ws.Range("F147").Value = 100112030
$ws.Range("G147").Value = "Poroto granado"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 38
$ws.Range("K147").Value = 25000
$ws.Range("L147").Value = 25000
$ws.Range("M147").Value = 25000
$ws.Range("N147").Value = "$/saco 25 kilos"
$ws.Range("O147").Value = "Provincia de Quillota"
$ws.Range("P147").Value = 1000
$ws.Range("Q147").Value = 25
$ws.Range("R147").Value = "Hortaliza"
